$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1152.4286
$ws.Range("I70").Value = 1008.4167
$ws.Range("J70").Value = 1344.4445
$ws.Range("K70").Value = 3025.2501
$ws.Range("L70").Value = 4033.3335
$ws.Range("M70").Value = -2755.2501
$ws.Range("N70").Value = -4573.333500000001
$ws.Range("H73").Value = 1152.4286
$ws.Range("I73").Value = 1008.4167
$ws.Range("J73").Value = 1344.4445
$ws.Range("K73").Value = 3025.2501
$ws.Range("L73").Value = 4033.3335
$ws.Range("M73").Value = -2089.2501
$ws.Range("N73").Value = -5905.333500000001
$ws.Range("H82").Value = 3321.0908
$ws.Range("I82").Value = 1154
$ws.Range("K82").Value = 3462
$ws.Range("M82").Value = -3056
$ws.Range("H85").Value = 3321.0908
$ws.Range("I85").Value = 1154
$ws.Range("K85").Value = 3462
$ws.Range("M85").Value = -2058
$ws.Range("H86").Value = 2868.7778
$ws.Range("I86").Value = 3999
$ws.Range("J86").Value = 2545.8572
$ws.Range("K86").Value = 3999
$ws.Range("L86").Value = 2545.8572
$ws.Range("M86").Value = -2876
$ws.Range("N86").Value = -4791.8572
$ws.Range("H89").Value = 2868.7778
$ws.Range("I89").Value = 3999
$ws.Range("J89").Value = 2545.8572
$ws.Range("K89").Value = 19995
$ws.Range("L89").Value = 12729.286
$ws.Range("M89").Value = -14379
$ws.Range("N89").Value = -23961.286
$ws.Range("H116").Value = 8146.6523
$ws.Range("I116").Value = 11061.154
$ws.Range("J116").Value = 4357.8
$ws.Range("K116").Value = 11061.154
$ws.Range("L116").Value = 4357.8
$ws.Range("M116").Value = -7619.154
$ws.Range("N116").Value = -11241.8
$ws.Range("H132").Value = 1729.3438
$ws.Range("I132").Value = 1611.7333
$ws.Range("K132").Value = 4835.199900000001
$ws.Range("M132").Value = -2305.199900000001
$ws.Range("H137").Value = 1061.7354
$ws.Range("I137").Value = 948.5625
$ws.Range("J137").Value = 1162.3334
$ws.Range("K137").Value = 2845.6875
$ws.Range("L137").Value = 3487.0002
$ws.Range("M137").Value = -295.6875
$ws.Range("N137").Value = -8587.0002
$ws.Range("H138").Value = 2984.4443
$ws.Range("I138").Value = 2427.5925
$ws.Range("J138").Value = 3402.0833
$ws.Range("K138").Value = 7282.7775
$ws.Range("L138").Value = 10206.2499
$ws.Range("M138").Value = -2142.7775
$ws.Range("N138").Value = -20486.2499
$ws.Range("H140").Value = 85526
$ws.Range("J140").Value = 85526
$ws.Range("L140").Value = 85526
$ws.Range("N140").Value = -95886

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1698.86
$ws.Range("I32").Value = 1698.86
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1698.86
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1411.86
$ws.Range("H74").Value = 2255.24
$ws.Range("I74").Value = 2104.5
$ws.Range("J74").Value = 2642.8572
$ws.Range("K74").Value = 2104.5
$ws.Range("L74").Value = 2642.8572
$ws.Range("M74").Value = -1230.5
$ws.Range("N74").Value = -4390.8572
$ws.Range("H77").Value = 2255.24
$ws.Range("I77").Value = 2104.5
$ws.Range("J77").Value = 2642.8572
$ws.Range("K77").Value = 10522.5
$ws.Range("L77").Value = 13214.286
$ws.Range("M77").Value = -6154.5
$ws.Range("N77").Value = -21950.286
$ws.Range("H88").Value = 4148.2144
$ws.Range("I88").Value = 1884.2858
$ws.Range("J88").Value = 4902.857
$ws.Range("K88").Value = 1884.2858
$ws.Range("L88").Value = 4902.857
$ws.Range("M88").Value = -1478.2858
$ws.Range("N88").Value = -5714.857
$ws.Range("H91").Value = 4148.2144
$ws.Range("I91").Value = 1884.2858
$ws.Range("J91").Value = 4902.857
$ws.Range("K91").Value = 1884.2858
$ws.Range("L91").Value = 4902.857
$ws.Range("M91").Value = -480.2858000000001
$ws.Range("N91").Value = -7710.857
$ws.Range("H138").Value = 61900
$ws.Range("J138").Value = 61900
$ws.Range("L138").Value = 61900
$ws.Range("N138").Value = -72180
$ws.Range("H139").Value = 60819.168
$ws.Range("J139").Value = 60819.168
$ws.Range("L139").Value = 60819.168
$ws.Range("N139").Value = -71099.16800000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6722.9355
$ws.Range("I134").Value = 1539.1111
$ws.Range("J134").Value = 13900.538
$ws.Range("K134").Value = 4617.3333
$ws.Range("L134").Value = 41701.614
$ws.Range("M134").Value = -2082.3333
$ws.Range("N134").Value = -46771.614
$ws.Range("H138").Value = 49495
$ws.Range("J138").Value = 49495
$ws.Range("L138").Value = 49495
$ws.Range("N138").Value = -59775
$ws.Range("H140").Value = 88325
$ws.Range("J140").Value = 88325
$ws.Range("L140").Value = 88325
$ws.Range("N140").Value = -98685

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 639.24
$ws.Range("I107").Value = 551.94116
$ws.Range("J107").Value = 824.75
$ws.Range("K107").Value = 551.94116
$ws.Range("L107").Value = 824.75
$ws.Range("M107").Value = 1368.05884
$ws.Range("N107").Value = -4664.75
$ws.Range("H132").Value = 1573.5106
$ws.Range("I132").Value = 1104.7567
$ws.Range("J132").Value = 3307.9
$ws.Range("K132").Value = 3314.2701
$ws.Range("L132").Value = 9923.700000000001
$ws.Range("M132").Value = -784.2700999999997
$ws.Range("N132").Value = -14983.7

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2526.1333
$ws.Range("I132").Value = 2123
$ws.Range("K132").Value = 6369
$ws.Range("M132").Value = -3839
$ws.Range("H136").Value = 10426.074
$ws.Range("J136").Value = 10426.074
$ws.Range("L136").Value = 31278.222
$ws.Range("N136").Value = -36378.222
$ws.Range("H138").Value = 61156.25
$ws.Range("J138").Value = 61156.25
$ws.Range("L138").Value = 61156.25
$ws.Range("N138").Value = -71436.25
$ws.Range("H140").Value = 99803
$ws.Range("J140").Value = 99803
$ws.Range("L140").Value = 99803
$ws.Range("N140").Value = -110163
$ws.Range("H141").Value = 66421.8
$ws.Range("J141").Value = 66421.8
$ws.Range("L141").Value = 66421.8
$ws.Range("N141").Value = -76781.8

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 91321.42999999999
$ws.Range("J133").Value = 91321.42999999999
$ws.Range("L133").Value = 91321.42999999999
$ws.Range("N133").Value = -96381.42999999999
$ws.Range("H136").Value = 2721.8235
$ws.Range("I136").Value = 1789.2703
$ws.Range("K136").Value = 5367.810899999999
$ws.Range("M136").Value = -2817.810899999999
$ws.Range("H139").Value = 79800
$ws.Range("J139").Value = 79800
$ws.Range("L139").Value = 79800
$ws.Range("N139").Value = -90080

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492
$ws.Range("H132").Value = 24591378
$ws.Range("I132").Value = 32609662
$ws.Range("J132").Value = 1971.7333
$ws.Range("K132").Value = 97828986
$ws.Range("L132").Value = 5915.199900000001
$ws.Range("M132").Value = -97826456
$ws.Range("N132").Value = -10975.1999
$ws.Range("H133").Value = 32700
$ws.Range("J133").Value = 32700
$ws.Range("L133").Value = 32700
$ws.Range("N133").Value = -42820
$ws.Range("H139").Value = 53942.855
$ws.Range("J139").Value = 53942.855
$ws.Range("L139").Value = 53942.855
$ws.Range("N139").Value = -64222.855
$ws.Range("H141").Value = 61857.145
$ws.Range("J141").Value = 64166.668
$ws.Range("L141").Value = 64166.668
$ws.Range("N141").Value = -74526.66800000001
